$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 5002
$ws.Range("I8").Value = 3
$ws.Range("K8").Value = 9
$ws.Range("M8").Value = 130
$ws.Range("H9").Value = 20.285715
$ws.Range("I9").Value = 19.5
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 19.5
$ws.Range("L9").Value = 25
$ws.Range("M9").Value = 149.5
$ws.Range("N9").Value = -363
$ws.Range("H19").Value = 300
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H28").Value = 1106.091
$ws.Range("I28").Value = 1151.1428
$ws.Range("J28").Value = 1027.25
$ws.Range("K28").Value = 1151.1428
$ws.Range("L28").Value = 1027.25
$ws.Range("M28").Value = -666.1428000000001
$ws.Range("N28").Value = -1997.25
$ws.Range("H38").Value = 6004
$ws.Range("I38").Value = 683.3333
$ws.Range("K38").Value = 2049.9999
$ws.Range("M38").Value = -1677.9999
$ws.Range("H98").Value = 739.0833
$ws.Range("I98").Value = 533
$ws.Range("K98").Value = 533
$ws.Range("M98").Value = 965
$ws.Range("H106").Value = 33949.777
$ws.Range("I106").Value = 35068.5
$ws.Range("K106").Value = 35068.5
$ws.Range("M106").Value = -34437.5
$ws.Range("H116").Value = 4821.5557
$ws.Range("I116").Value = 4663.4287
$ws.Range("K116").Value = 4663.4287
$ws.Range("M116").Value = -1221.4287
$ws.Range("H122").Value = 739.0833
$ws.Range("I122").Value = 533
$ws.Range("K122").Value = 1599
$ws.Range("M122").Value = 851

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8823.518
$ws.Range("I32").Value = 4603.696
$ws.Range("K32").Value = 4603.696
$ws.Range("M32").Value = -4316.696
$ws.Range("H45").Value = 1740.5
$ws.Range("I45").Value = 1712
$ws.Range("K45").Value = 1712
$ws.Range("M45").Value = -1335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3364.5715
$ws.Range("I99").Value = 3211.4666
$ws.Range("J99").Value = 3747.3333
$ws.Range("K99").Value = 3211.4666
$ws.Range("L99").Value = 3747.3333
$ws.Range("M99").Value = -1713.4666
$ws.Range("N99").Value = -6743.3333
$ws.Range("H107").Value = 615.44446
$ws.Range("I107").Value = 612.72
$ws.Range("K107").Value = 612.72
$ws.Range("M107").Value = 1307.28

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3659.7
$ws.Range("I31").Value = 1592.2632
$ws.Range("J31").Value = 7230.727
$ws.Range("K31").Value = 1592.2632
$ws.Range("L31").Value = 7230.727
$ws.Range("M31").Value = -1297.2632
$ws.Range("N31").Value = -7820.727
$ws.Range("H34").Value = 3659.7
$ws.Range("I34").Value = 1592.2632
$ws.Range("J34").Value = 7230.727
$ws.Range("K34").Value = 1592.2632
$ws.Range("L34").Value = 7230.727
$ws.Range("M34").Value = -1390.2632
$ws.Range("N34").Value = -7634.727
$ws.Range("H86").Value = 10499.25
$ws.Range("I86").Value = 7998
$ws.Range("J86").Value = 11333
$ws.Range("K86").Value = 7998
$ws.Range("L86").Value = 11333
$ws.Range("M86").Value = -6875
$ws.Range("N86").Value = -13579
$ws.Range("H89").Value = 10499.25
$ws.Range("I89").Value = 7998
$ws.Range("J89").Value = 11333
$ws.Range("K89").Value = 39990
$ws.Range("L89").Value = 56665
$ws.Range("M89").Value = -34374
$ws.Range("N89").Value = -67897
$ws.Range("H99").Value = 12884.969
$ws.Range("I99").Value = 9529.267
$ws.Range("J99").Value = 15845.883
$ws.Range("K99").Value = 9529.267
$ws.Range("L99").Value = 15845.883
$ws.Range("M99").Value = -8031.267
$ws.Range("N99").Value = -18841.883
$ws.Range("H107").Value = 2150.7666
$ws.Range("I107").Value = 2332.6086
$ws.Range("K107").Value = 2332.6086
$ws.Range("M107").Value = -412.6086
$ws.Range("H122").Value = 3264.4443
$ws.Range("I122").Value = 3409.7856
$ws.Range("K122").Value = 10229.3568
$ws.Range("M122").Value = -7779.356800000001
$ws.Range("H126").Value = 12884.969
$ws.Range("I126").Value = 9529.267
$ws.Range("J126").Value = 15845.883
$ws.Range("K126").Value = 28587.801
$ws.Range("L126").Value = 47537.649
$ws.Range("M126").Value = -26117.801
$ws.Range("N126").Value = -52477.649

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 452
$ws.Range("I5").Value = 520.625
$ws.Range("J5").Value = 412.7857
$ws.Range("K5").Value = 1561.875
$ws.Range("L5").Value = 1238.3571
$ws.Range("M5").Value = -1449.875
$ws.Range("N5").Value = -1462.3571
$ws.Range("H8").Value = 188.28572
$ws.Range("I8").Value = 188.28572
$ws.Range("K8").Value = 564.85716
$ws.Range("M8").Value = -425.85716
$ws.Range("H68").Value = 1262.6
$ws.Range("I68").Value = 1153.25
$ws.Range("K68").Value = 3459.75
$ws.Range("M68").Value = -2648.75
$ws.Range("H71").Value = 1262.6
$ws.Range("I71").Value = 1153.25
$ws.Range("K71").Value = 10379.25
$ws.Range("M71").Value = -6323.25
$ws.Range("H92").Value = 528.1852
$ws.Range("J92").Value = 561.53845
$ws.Range("L92").Value = 1684.61535
$ws.Range("N92").Value = -4180.61535
$ws.Range("H122").Value = 328.14285
$ws.Range("I122").Value = 356.66666
$ws.Range("J122").Value = 306.75
$ws.Range("K122").Value = 3209.99994
$ws.Range("L122").Value = 2760.75
$ws.Range("M122").Value = -759.9999399999997
$ws.Range("N122").Value = -7660.75
$ws.Range("H135").Value = 452
$ws.Range("I135").Value = 520.625
$ws.Range("J135").Value = 412.7857
$ws.Range("K135").Value = 4685.625
$ws.Range("L135").Value = 3715.0713
$ws.Range("M135").Value = -2150.625
$ws.Range("N135").Value = -8785.0713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5000
$ws.Range("J70").Value = 5000
$ws.Range("L70").Value = 5000
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 5000
$ws.Range("J73").Value = 5000
$ws.Range("L73").Value = 5000
$ws.Range("N73").Value = -6872
$ws.Range("H97").Value = 919.1818
$ws.Range("I97").Value = 748
$ws.Range("J97").Value = 1375.6666
$ws.Range("K97").Value = 748
$ws.Range("L97").Value = 1375.6666
$ws.Range("M97").Value = -252
$ws.Range("N97").Value = -2367.6666
$ws.Range("H102").Value = 4812.364
$ws.Range("I102").Value = 4487
$ws.Range("J102").Value = 4998.2856
$ws.Range("K102").Value = 4487
$ws.Range("L102").Value = 4998.2856
$ws.Range("M102").Value = -2865
$ws.Range("N102").Value = -8242.285599999999
$ws.Range("H132").Value = 5366.222
$ws.Range("I132").Value = 2005.5
$ws.Range("J132").Value = 8054.8
$ws.Range("K132").Value = 6016.5
$ws.Range("L132").Value = 24164.4
$ws.Range("M132").Value = -3486.5
$ws.Range("N132").Value = -29224.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 317.625
$ws.Range("I55").Value = 319.6
$ws.Range("K55").Value = 319.6
$ws.Range("M55").Value = -146.6
$ws.Range("H82").Value = 3196.76
$ws.Range("I82").Value = 3436.2778
$ws.Range("J82").Value = 2580.8572
$ws.Range("K82").Value = 3436.2778
$ws.Range("L82").Value = 2580.8572
$ws.Range("M82").Value = -3075.2778
$ws.Range("N82").Value = -3302.8572
$ws.Range("H85").Value = 3196.76
$ws.Range("I85").Value = 3436.2778
$ws.Range("J85").Value = 2580.8572
$ws.Range("K85").Value = 3436.2778
$ws.Range("L85").Value = 2580.8572
$ws.Range("M85").Value = -2188.2778
$ws.Range("N85").Value = -5076.8572
$ws.Range("H122").Value = 8499.75
$ws.Range("I122").Value = 8499.75
$ws.Range("K122").Value = 25499.25
$ws.Range("M122").Value = -23049.25
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("N131").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6451.2
$ws.Range("I62").Value = 5534.625
$ws.Range("J62").Value = 7498.7144
$ws.Range("K62").Value = 5534.625
$ws.Range("L62").Value = 7498.7144
$ws.Range("M62").Value = -4910.625
$ws.Range("N62").Value = -8746.714400000001
$ws.Range("H65").Value = 6451.2
$ws.Range("I65").Value = 5534.625
$ws.Range("J65").Value = 7498.7144
$ws.Range("K65").Value = 27673.125
$ws.Range("L65").Value = 37493.572
$ws.Range("M65").Value = -24553.125
$ws.Range("N65").Value = -43733.572
$ws.Range("H100").Value = 2034.1
$ws.Range("I100").Value = 1732.8334
$ws.Range("J100").Value = 2486
$ws.Range("K100").Value = 3465.6668
$ws.Range("L100").Value = 4972
$ws.Range("M100").Value = -2924.6668
$ws.Range("N100").Value = -6054
$ws.Range("H113").Value = 468.33334
$ws.Range("I113").Value = 593.4167
$ws.Range("K113").Value = 1780.2501
$ws.Range("M113").Value = 389.7499
$ws.Range("I122").Value = 10499.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 31498.5
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -29048.5
$ws.Range("H126").Value = 1540.2106
$ws.Range("I126").Value = 1078.6875
$ws.Range("J126").Value = 4001.6667
$ws.Range("K126").Value = 3236.0625
$ws.Range("L126").Value = 12005.0001
$ws.Range("M126").Value = -766.0625
$ws.Range("N126").Value = -16945.0001
$ws.Range("H132").Value = 1469.8096
$ws.Range("I132").Value = 1329.125
$ws.Range("K132").Value = 3987.375
$ws.Range("M132").Value = -1457.375
